# NIT-9014769947.xlsx update
# - Updates VALOR MORA / Cant. Trabajadores / Cant. Periodos header values
# - Re-orders/extends the periodo rows for the existing worker (DAVID ANTONIO DIAZ SALAS)
# - Appends 13 new worker rows (part 1 of new estado de cuenta)
# - Moves the signature block down below the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Preserve formatting for the two "template" row styles we need:
#    - the special bottom-border style currently on row 19 (last data row)
#    - the plain middle style currently on row 18
#    Do this BEFORE we touch rows 19/24/25 so the source formatting is intact.
# ------------------------------------------------------------------

# Move the special "last row" look from row 19 down to row 32 (new last data row)
$ws.Range("B19:J19").Copy($ws.Range("B32:J32"))

# Move the signature block (rows 24-25) down to rows 37-38 before we overwrite 24-25
$ws.Range("B24:C24").Copy($ws.Range("B37:C37"))
$ws.Range("H24:J24").Copy($ws.Range("H37:J37"))
$ws.Range("B25:C25").Copy($ws.Range("B38:C38"))
$ws.Range("H25:J25").Copy($ws.Range("H38:J38"))

# Free up the old merges on rows 24 and 25 - they become plain data rows
$ws.Range("B24:C24").UnMerge()
$ws.Range("H24:J24").UnMerge()
$ws.Range("B25:C25").UnMerge()
$ws.Range("H25:J25").UnMerge()

# Copy() skips blank source cells instead of clearing the destination, so the
# old signature text lingering in H24:J24 / H25:J25 needs an explicit wipe.
$ws.Range("H24:J24").ClearContents()
$ws.Range("H25:J25").ClearContents()

# Stamp the plain middle style (from row 18) onto row 19 and the 12 new rows (20-31)
$ws.Range("B18:J18").Copy($ws.Range("B19:J19"))
$ws.Range("B18:J18").Copy($ws.Range("B20:J20"))
$ws.Range("B18:J18").Copy($ws.Range("B21:J21"))
$ws.Range("B18:J18").Copy($ws.Range("B22:J22"))
$ws.Range("B18:J18").Copy($ws.Range("B23:J23"))
$ws.Range("B18:J18").Copy($ws.Range("B24:J24"))
$ws.Range("B18:J18").Copy($ws.Range("B25:J25"))
$ws.Range("B18:J18").Copy($ws.Range("B26:J26"))
$ws.Range("B18:J18").Copy($ws.Range("B27:J27"))
$ws.Range("B18:J18").Copy($ws.Range("B28:J28"))
$ws.Range("B18:J18").Copy($ws.Range("B29:J29"))
$ws.Range("B18:J18").Copy($ws.Range("B30:J30"))
$ws.Range("B18:J18").Copy($ws.Range("B31:J31"))

# ------------------------------------------------------------------
# 2. Header values
# ------------------------------------------------------------------
$ws.Range("E11").Value = 874887
$ws.Range("C13").Value = 14
$ws.Range("F13").Value = 5

# ------------------------------------------------------------------
# 3. DAVID ANTONIO DIAZ SALAS periodo rows (16-19), now chronological
# ------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1003080810"
$ws.Range("D16").Value = "DAVID ANTONIO DIAZ SALAS"
$ws.Range("E16").Value = "2212"
$ws.Range("F16").Value = 14667
$ws.Range("G16").Value = 1423500

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1003080810"
$ws.Range("D17").Value = "DAVID ANTONIO DIAZ SALAS"
$ws.Range("E17").Value = "2301"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1003080810"
$ws.Range("D18").Value = "DAVID ANTONIO DIAZ SALAS"
$ws.Range("E18").Value = "2302"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1003080810"
$ws.Range("D19").Value = "DAVID ANTONIO DIAZ SALAS"
$ws.Range("E19").Value = "2303"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1423500

# ------------------------------------------------------------------
# 4. New workers, all periodo 2508 (rows 20-32)
# ------------------------------------------------------------------
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "80882819"
$ws.Range("D20").Value = "FRANCISCO DIAZ CARMONA"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "22803653"
$ws.Range("D21").Value = "KATIA VEGA GOMEZ"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "45559163"
$ws.Range("D22").Value = "CLAUDIA DEL CARMEN DEL CASTILLO RAMIREZ"
$ws.Range("E22").Value = "2508"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1103103788"
$ws.Range("D23").Value = "JEISON ARLEY ARRIETA ORTEGA"
$ws.Range("E23").Value = "2508"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "9146478"
$ws.Range("D24").Value = "CARLOS ANDRES CRISMATH CASTIILO"
$ws.Range("E24").Value = "2508"
$ws.Range("F24").Value = 56940
$ws.Range("G24").Value = 1423500

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1143388520"
$ws.Range("D25").Value = "ALFREDO GUZMAN PALACIO"
$ws.Range("E25").Value = "2508"
$ws.Range("F25").Value = 56940
$ws.Range("G25").Value = 1423500

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "5165754"
$ws.Range("D26").Value = "ADRIAN GREGORIO MERCADO VERGARA"
$ws.Range("E26").Value = "2508"
$ws.Range("F26").Value = 56940
$ws.Range("G26").Value = 1423500

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1237441882"
$ws.Range("D27").Value = "EDUARDO LUIS OSPINO TORRES"
$ws.Range("E27").Value = "2508"
$ws.Range("F27").Value = 56940
$ws.Range("G27").Value = 1423500

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1002242188"
$ws.Range("D28").Value = "ISAAC DAVID BELLO RODRIGUEZ"
$ws.Range("E28").Value = "2508"
$ws.Range("F28").Value = 56940
$ws.Range("G28").Value = 1423500

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "98590856"
$ws.Range("D29").Value = "OSSMAN ELIECER BUITRAGO CASTAÝO"
$ws.Range("E29").Value = "2508"
$ws.Range("F29").Value = 56940
$ws.Range("G29").Value = 1423500

$ws.Range("B30").Value = "PPT"
$ws.Range("C30").Value = "6687758"
$ws.Range("D30").Value = "RICARDO JESUS FERNANDEZ SILVA"
$ws.Range("E30").Value = "2508"
$ws.Range("F30").Value = 56940
$ws.Range("G30").Value = 1423500

$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "9290558"
$ws.Range("D31").Value = "JOSE ANGEL BERMEJO CARDONA"
$ws.Range("E31").Value = "2508"
$ws.Range("F31").Value = 56940
$ws.Range("G31").Value = 1423500

$ws.Range("B32").Value = "PPT"
$ws.Range("C32").Value = "5646957"
$ws.Range("D32").Value = "FRANCO RAFAEL OSPINO TORRES"
$ws.Range("E32").Value = "2508"
$ws.Range("F32").Value = 56940
$ws.Range("G32").Value = 1423500

# ------------------------------------------------------------------
# 5. Column D needs to be wide enough for the longest new name
# ------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 43.08984375
